$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (price) cells from Excel auto-numeric conversion,
# since the source data stores these as literal text (inline strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.373.62"
$ws.Range("E2").Value = "  +2.84%  "

$ws.Range("D3").Value = "2.305.79"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "311.41"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").Value = "102.82"
$ws.Range("E6").Value = "  +6.50%  "

$ws.Range("E7").Value = "  +1.56%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +7.68%  "

$ws.Range("D10").Value = "35.88"
$ws.Range("E10").Value = "  +2.59%  "

$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +3.05%  "

$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("D13").Value = "6.99"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").Value = "2.664.29"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").Value = "15.04"
$ws.Range("E15").Value = "  +2.54%  "

$ws.Range("D16").Value = "2.315.88"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").Value = "0.809"
$ws.Range("E17").Value = "  +2.24%  "

$ws.Range("D18").Value = "43.316.24"
$ws.Range("E18").Value = "  +3.01%  "

$ws.Range("D19").Value = "12.34"
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").Value = "0.0₃0934"
$ws.Range("E20").Value = "  +3.31%  "

$ws.Range("E21").Value = "  +3.15%  "

$ws.Range("D22").Value = "68.13"
$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("D23").Value = "241.37"
$ws.Range("E23").Value = "  +1.83%  "

$ws.Range("D24").Value = "2.62"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("E25").Value = "  +2.76%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "24.64"
$ws.Range("E27").Value = "  +4.90%  "

$ws.Range("E28").Value = "  +8.57%  "

$ws.Range("D29").Value = "37.07"
$ws.Range("E29").Value = "  -1.45%  "

$ws.Range("D30").Value = "9.65"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").Value = "169.04"
$ws.Range("E31").Value = "  +3.90%  "

$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  +6.25%  "

$ws.Range("D35").Value = "0.0743"
$ws.Range("E35").Value = "  +0.89%  "

$ws.Range("D36").Value = "17.64"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -3.29%  "

$ws.Range("D38").Value = "1.89"
$ws.Range("E38").Value = "  +3.80%  "

$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("E40").Value = "  +1.64%  "

$ws.Range("D41").Value = "4.38"
$ws.Range("E41").Value = "  +7.94%  "

$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").Value = "19.58"
$ws.Range("E43").Value = "  +2.87%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0290"
$ws.Range("E44").Value = "  +3.31%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.972.80"
$ws.Range("E45").Value = "  +1.18%  "

$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +2.31%  "

$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").Value = "55.41"
$ws.Range("E48").Value = "  +2.56%  "

$ws.Range("D49").Value = "2.92"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("E50").Value = "  +7.90%  "

$ws.Range("D51").Value = "2.533.72"

# Restore the default style on column D so only the values changed
# (matches the original workbook, which had no explicit numeric format).
$ws.Range("D2:D51").Style = "Normal"
